# Edit: for each year block of 4 quarterly rows (A,B,C,D) starting at row 2,
# swap the B-row and C-row contents (columns A:E), then delete columns F:G
# entirely (the "钢材产销率" / "钢材销售量" quarter-only columns were dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (B-row, C-row) to swap for each year, columns A:E.
for ($r = 3; $r -le 79; $r += 4) {
    $r2 = $r + 1
    $rowB = $ws.Range("A$r`:E$r").Value()
    $rowC = $ws.Range("A$r2`:E$r2").Value()
    $ws.Range("A$r`:E$r").Value = $rowC
    $ws.Range("A$r2`:E$r2").Value = $rowB
}

# Drop the now-unused F (钢材产销率) and G (钢材销售量) columns.
$ws.Range("F:G").Delete()
